# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" columns (E/F)
# for the first two data rows of the zh-cn and de-de sheets, stamps the
# "Latest Handback DateTime" (column G), and flips the Overview/per-locale
# Status text from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

function Get-LinkTarget($sheet, $addr) {
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            return $h.Address
        }
    }
    return $null
}

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the Status column (B/C) mirrors the same shared text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# Helper that fills in one locale sheet (zh-cn / de-de).
# ---------------------------------------------------------------------
function Update-LocaleSheet($sheetName, $handbackFileName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (B) for the two real rows.
    $ws.Range("B2").Value = $statusHandedBack
    $ws.Range("B3").Value = $statusHandedBack

    # Targets for the new hyperlinks reuse the exact same links already
    # used by column A (source file) and column C (handoff xlf) for each row.
    $targetA2 = Get-LinkTarget $ws '$A$2'
    $targetC2 = Get-LinkTarget $ws '$C$2'
    $targetA3 = Get-LinkTarget $ws '$A$3'
    $targetC3 = Get-LinkTarget $ws '$C$3'

    $sourceFileDisplay = $ws.Range("A2").Text
    $handoffFileDisplay = $handbackFileName

    # Row 2: Latest Target File (E2) / Latest Handback File (F2).
    $ws.Range("E2").Value = $sourceFileDisplay
    $ws.Hyperlinks.Add($ws.Range("E2"), $targetA2, "", "", $sourceFileDisplay) | Out-Null

    $ws.Range("F2").Value = $handoffFileDisplay
    $ws.Hyperlinks.Add($ws.Range("F2"), $targetC2, "", "", $handoffFileDisplay) | Out-Null

    $ws.Range("G2").Value = $handbackDateTime

    # Row 3: Latest Target File (E3) / Latest Handback File (F3).
    $ws.Range("E3").Value = $sourceFileDisplay
    $ws.Hyperlinks.Add($ws.Range("E3"), $targetA3, "", "", $sourceFileDisplay) | Out-Null

    $ws.Range("F3").Value = $handoffFileDisplay
    $ws.Hyperlinks.Add($ws.Range("F3"), $targetC3, "", "", $handoffFileDisplay) | Out-Null

    $ws.Range("G3").Value = $handbackDateTime
}

Update-LocaleSheet "zh-cn" "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf" "2016-03-01 09:39:07"
Update-LocaleSheet "de-de" "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf" "2016-03-01 09:39:25"
